# Fix multiple genes in alteration table
# Adds a "join_at" column to the CDS!Tabella2 table and a new row for the
# second ORF1ab (join_at) segment.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)   # "CDS" sheet
$lo = $ws.ListObjects.Item(1)  # "Tabella2"

# --- add the new "join_at" table column -----------------------------------
$joinCol = $lo.ListColumns.Add()
$joinCol.Range.Cells.Item(1, 1).Value = "join_at"

# Give the new column's (currently empty) data cells the same formatting as
# the existing "to" column (s="2": Arial Unicode MS 10pt, vertical centred).
$ws.Range("D2").Copy() | Out-Null
$ws.Range("E2:E12").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null
$excel.CutCopyMode = $false

# --- add the new data row (second ORF1ab / join_at segment) ---------------
$newRow = $lo.ListRows.Add()

# Inherit formatting from the row above before filling in values.
$ws.Range("A12:E12").Copy() | Out-Null
$ws.Range("A13:E13").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("A13").Value = "ORF1ab"
$ws.Range("B13").Value = 43740578
$ws.Range("C13").Value = 266
$ws.Range("D13").Value = 21555
$ws.Range("E13").Value = 13468

# --- match the author's final selection state ------------------------------
$ws.Range("E13").Select() | Out-Null
